$d = $word.ActiveDocument

# The "References" section ends with three hyperlink paragraphs followed by
# a single empty paragraph right before the sectPr. We need to insert a new
# paragraph containing a hyperlink to the StackOverflow question, directly
# after the last existing hyperlink (CodeProject) and before that trailing
# empty paragraph.

# Locate the last hyperlink paragraph (the CodeProject reference) - it is
# the paragraph immediately preceding the final, empty paragraph.
$lastParaIndex = $d.Paragraphs.Count
$trailingEmptyPara = $d.Paragraphs.Item($lastParaIndex)
$lastLinkPara = $d.Paragraphs.Item($lastParaIndex - 1)

# Add the new hyperlink right at the end of the CodeProject paragraph; this
# folds cleanly into a brand-new paragraph (Word splits on the paragraph
# mark) without leaving stray empty runs behind.
$insertionPoint = $d.Range($lastLinkPara.Range.End, $lastLinkPara.Range.End)
$d.Hyperlinks.Add($insertionPoint, `
    "http://stackoverflow.com/questions/4718725/best-practice-to-avoid-invalidoperationexception-collection-was-modified")

# Restore the trailing empty paragraph that InsertParagraphAfter/Hyperlinks.Add
# consumed, by splitting a fresh paragraph mark off the end of the document.
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfDoc = $d.Range($newLastPara.Range.End, $newLastPara.Range.End)
$endOfDoc.Text = [char]13
